$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived statistics for the 16-row x 14-column block G2:T17
$arr = New-Object 'object[,]' 16,14
$arr[0,0] = 2.374260666666666
$arr[0,1] = 7.122781999999999
$arr[0,2] = 0.0276017086472712
$arr[0,3] = 0.0276017086472712
$arr[0,4] = 3.0
$arr[0,5] = 1.0
$arr[0,6] = 0.312815
$arr[0,7] = 0.938445
$arr[0,8] = 0.0082131704949067
$arr[0,9] = 0.0082131704949067
$arr[0,10] = 0.7427043504433333
$arr[0,11] = 6.684339153989999
$arr[0,12] = 0.0002266975390707789
$arr[0,13] = 0.0002266975390707789
$arr[1,0] = 2.374260666666666
$arr[1,1] = 7.122781999999999
$arr[1,2] = 0.0276017086472712
$arr[1,3] = 0.0276017086472712
$arr[1,4] = 3.0
$arr[1,5] = 1.0
$arr[1,6] = 24.84824866666667
$arr[1,7] = 74.544746
$arr[1,8] = 0.6524076620340182
$arr[1,9] = 0.6524076620340182
$arr[1,10] = 58.9962194448191
$arr[1,11] = 530.9659750033719
$arr[1,12] = 0.01800756620671035
$arr[1,13] = 0.01800756620671035
$arr[2,0] = 2.374260666666666
$arr[2,1] = 7.122781999999999
$arr[2,2] = 0.0276017086472712
$arr[2,3] = 0.0276017086472712
$arr[2,4] = 3.0
$arr[2,5] = 1.0
$arr[2,6] = 12.866992
$arr[2,7] = 38.600976
$arr[2,8] = 0.3378316226926476
$arr[2,9] = 0.3378316226926476
$arr[2,10] = 30.54959300391467
$arr[2,11] = 274.946337035232
$arr[2,12] = 0.009324730021397313
$arr[2,13] = 0.009324730021397313
$arr[3,0] = 2.374260666666666
$arr[3,1] = 7.122781999999999
$arr[3,2] = 0.0276017086472712
$arr[3,3] = 0.0276017086472712
$arr[3,4] = 3.0
$arr[3,5] = 1.0
$arr[3,6] = 0.05894133333333334
$arr[3,7] = 0.176824
$arr[3,8] = 0.001547544778427486
$arr[3,9] = 0.001547544778427486
$arr[3,10] = 0.1399420893742222
$arr[3,11] = 1.259478804368
$arr[3,12] = 0.00004271488009276134
$arr[3,13] = 0.00004271488009276134
$arr[4,0] = 76.92488366666667
$arr[4,1] = 230.774651
$arr[4,2] = 0.8942818522422411
$arr[4,3] = 0.8942818522422411
$arr[4,4] = 3.0
$arr[4,5] = 1.0
$arr[4,6] = 0.312815
$arr[4,7] = 0.938445
$arr[4,8] = 0.0082131704949067
$arr[4,9] = 0.0082131704949067
$arr[4,10] = 24.06325748418834
$arr[4,11] = 216.569317357695
$arr[4,12] = 0.007344889322966487
$arr[4,13] = 0.007344889322966487
$arr[5,0] = 76.92488366666667
$arr[5,1] = 230.774651
$arr[5,2] = 0.8942818522422411
$arr[5,3] = 0.8942818522422411
$arr[5,4] = 3.0
$arr[5,5] = 1.0
$arr[5,6] = 24.84824866666667
$arr[5,7] = 74.544746
$arr[5,8] = 0.6524076620340182
$arr[5,9] = 0.6524076620340182
$arr[5,10] = 1911.448638003739
$arr[5,11] = 17203.03774203365
$arr[5,12] = 0.5834363324208118
$arr[5,13] = 0.5834363324208118
$arr[6,0] = 76.92488366666667
$arr[6,1] = 230.774651
$arr[6,2] = 0.8942818522422411
$arr[6,3] = 0.8942818522422411
$arr[6,4] = 3.0
$arr[6,5] = 1.0
$arr[6,6] = 12.866992
$arr[6,7] = 38.600976
$arr[6,8] = 0.3378316226926476
$arr[6,9] = 0.3378316226926476
$arr[6,10] = 989.7918627399308
$arr[6,11] = 8908.126764659377
$arr[6,12] = 0.3021166892875828
$arr[6,13] = 0.3021166892875828
$arr[7,0] = 76.92488366666667
$arr[7,1] = 230.774651
$arr[7,2] = 0.8942818522422411
$arr[7,3] = 0.8942818522422411
$arr[7,4] = 3.0
$arr[7,5] = 1.0
$arr[7,6] = 0.05894133333333334
$arr[7,7] = 0.176824
$arr[7,8] = 0.001547544778427486
$arr[7,9] = 0.001547544778427486
$arr[7,10] = 4.534055209824889
$arr[7,11] = 40.806496888424
$arr[7,12] = 0.001383941210879941
$arr[7,13] = 0.001383941210879941
$arr[8,0] = 6.625048
$arr[8,1] = 19.875144
$arr[8,2] = 0.07701877356495823
$arr[8,3] = 0.07701877356495825
$arr[8,4] = 3.0
$arr[8,5] = 1.0
$arr[8,6] = 0.312815
$arr[8,7] = 0.938445
$arr[8,8] = 0.0082131704949067
$arr[8,9] = 0.0082131704949067
$arr[8,10] = 2.07241439012
$arr[8,11] = 18.65172951108
$arr[8,12] = 0.000632568318597615
$arr[8,13] = 0.0006325683185976151
$arr[9,0] = 6.625048
$arr[9,1] = 19.875144
$arr[9,2] = 0.07701877356495823
$arr[9,3] = 0.07701877356495825
$arr[9,4] = 3.0
$arr[9,5] = 1.0
$arr[9,6] = 24.84824866666667
$arr[9,7] = 74.544746
$arr[9,8] = 0.6524076620340182
$arr[9,9] = 0.6524076620340182
$arr[9,10] = 164.6208401326027
$arr[9,11] = 1481.587561193424
$arr[9,12] = 0.05024763799424185
$arr[9,13] = 0.05024763799424185
$arr[10,0] = 6.625048
$arr[10,1] = 19.875144
$arr[10,2] = 0.07701877356495823
$arr[10,3] = 0.07701877356495825
$arr[10,4] = 3.0
$arr[10,5] = 1.0
$arr[10,6] = 12.866992
$arr[10,7] = 38.600976
$arr[10,8] = 0.3378316226926476
$arr[10,9] = 0.3378316226926476
$arr[10,10] = 85.244439615616
$arr[10,11] = 767.199956540544
$arr[10,12] = 0.02601937725124743
$arr[10,13] = 0.02601937725124744
$arr[11,0] = 6.625048
$arr[11,1] = 19.875144
$arr[11,2] = 0.07701877356495823
$arr[11,3] = 0.07701877356495825
$arr[11,4] = 3.0
$arr[11,5] = 1.0
$arr[11,6] = 0.05894133333333334
$arr[11,7] = 0.176824
$arr[11,8] = 0.001547544778427486
$arr[11,9] = 0.001547544778427486
$arr[11,10] = 0.3904891625173333
$arr[11,11] = 3.514402462656
$arr[11,12] = 0.00011919000087134
$arr[11,13] = 0.00011919000087134
$arr[12,0] = 0.09441966666666667
$arr[12,1] = 0.283259
$arr[12,2] = 0.001097665545529457
$arr[12,3] = 0.001097665545529457
$arr[12,4] = 3.0
$arr[12,5] = 1.0
$arr[12,6] = 0.312815
$arr[12,7] = 0.938445
$arr[12,8] = 0.0082131704949067
$arr[12,9] = 0.0082131704949067
$arr[12,10] = 0.02953588802833334
$arr[12,11] = 0.265822992255
$arr[12,12] = 0.000009015314271818199
$arr[12,13] = 0.000009015314271818199
$arr[13,0] = 0.09441966666666667
$arr[13,1] = 0.283259
$arr[13,2] = 0.001097665545529457
$arr[13,3] = 0.001097665545529457
$arr[13,4] = 3.0
$arr[13,5] = 1.0
$arr[13,6] = 24.84824866666667
$arr[13,7] = 74.544746
$arr[13,8] = 0.6524076620340182
$arr[13,9] = 0.6524076620340182
$arr[13,10] = 2.346163356357111
$arr[13,11] = 21.115470207214
$arr[13,12] = 0.0007161254122541679
$arr[13,13] = 0.0007161254122541679
$arr[14,0] = 0.09441966666666667
$arr[14,1] = 0.283259
$arr[14,2] = 0.001097665545529457
$arr[14,3] = 0.001097665545529457
$arr[14,4] = 3.0
$arr[14,5] = 1.0
$arr[14,6] = 12.866992
$arr[14,7] = 38.600976
$arr[14,8] = 0.3378316226926476
$arr[14,9] = 0.3378316226926476
$arr[14,10] = 1.214897095642667
$arr[14,11] = 10.934073860784
$arr[14,12] = 0.0003708261324200265
$arr[14,13] = 0.0003708261324200265
$arr[15,0] = 0.09441966666666667
$arr[15,1] = 0.283259
$arr[15,2] = 0.001097665545529457
$arr[15,3] = 0.001097665545529457
$arr[15,4] = 3.0
$arr[15,5] = 1.0
$arr[15,6] = 0.05894133333333334
$arr[15,7] = 0.176824
$arr[15,8] = 0.001547544778427486
$arr[15,9] = 0.001547544778427486
$arr[15,10] = 0.005565221046222223
$arr[15,11] = 0.050086989416
$arr[15,12] = 0.000001698686583443869
$arr[15,13] = 0.000001698686583443869

$ws.Range("G2:T17").Value = $arr
